$d = $word.ActiveDocument

# --- Change 1: remove the leading "scratch notes" paragraphs -------------
# The document used to start with a date stamp, a blank line, and three
# dash-bullet TODO lines. Those five paragraphs are removed entirely
# (including their paragraph marks), leaving the two blank paragraphs that
# originally followed them.
$firstPara = $d.Paragraphs.Item(1)
$fifthPara = $d.Paragraphs.Item(5)
$deleteRange = $d.Range($firstPara.Range.Start, $fifthPara.Range.End)
$deleteRange.Delete()

# --- Change 2: tag the "-add time (1 minute)" run with a rendered page
#     break marker (<w:lastRenderedPageBreak/>) as the first child of its
#     run, right before the text. -----------------------------------------
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "-add time (1 minute)`r") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $full = $r.WordOpenXML
    if ($full -match '(<w:p\b[^>]*>.*?</w:p>)') {
        $paraXml = $matches[1]
        # insert the break marker as the very first child of the first run
        $newParaXml = $paraXml -replace '<w:r>', '<w:r><w:lastRenderedPageBreak/>', 1
        $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
               '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $r.InsertXML($pkg)
    }
}
